$d = $word.ActiveDocument

$d.Content.Find.Execute("Start time: 2017-12-27 18:30:19", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Start time: 2018-01-31 12:32:50", 2)

$d.Content.Find.Execute("End time: 2017-12-27 18:30:28", $true, $false, $false, $false, $false,
                         $true, 1, $false, "End time: 2018-01-31 12:32:59", 2)

$d.Content.Find.Execute("Duration: 9.39 secs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Duration: 8.88 secs", 2)
